$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Domeinen")

# --- New "wetgeving" (legislation) columns D & E, added while prepping the Archi import ---
# Header
$ws.Range("D1").Value = "wetgeving"

# HAP (Horeca en ambachtelijke productie) - row 17
$ws.Range("D17").Value = "BWBR0002458"

# AT (Alcohol en tabak) - row 3
$ws.Range("D3").Value = "BWBR0002458, BWBR0004302"
$ws.Range("D3").Font.Name = "Arial"

# DP (Dierproeven) - row 10
$ws.Range("D10").Value = "BWBR0003081"
$ws.Range("D10").Font.Name = "Arial"
$ws.Range("E10").Value = "Wet op de dierproeven"
$ws.Range("E10").Font.Name = "Arial"

# HAP (Horeca en ambachtelijke productie) - row 17
$ws.Range("E17").Value = "Drank- en Horecawet"

# AT (Alcohol en tabak) - row 3
$ws.Range("E3").Value = "Drank- en Horecawet; Tabakswet"
$ws.Range("E3").Font.Name = "Arial"

# DGM (Diergeneesmiddelen) - row 8
$ws.Range("D8").Value = "BWBR0003818"
$ws.Range("D8").Font.Name = "Arial"
$ws.Range("E8").Value = "Diergeneesmiddelenwet"
$ws.Range("E8").Font.Name = "Arial"

# DW (Dierenwelzijn) - row 7
$ws.Range("D7").Value = "BWBR0005662"
$ws.Range("D7").Font.Name = "Arial"
$ws.Range("E7").Value = "Gezondheids- en welzijnswet voor dieren"
$ws.Range("E7").Font.Name = "Arial"

# Widen columns C & D to fit the new legislation text
$ws.Columns("C").ColumnWidth = 35
$ws.Columns("D").ColumnWidth = 43.6

# Make "Domeinen" the active/visible sheet again (was "Divisies" before)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A23").Select()
